$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: paragraph about the IGBT cycloconverter model - insert a new
# sentence ("As can be observed ... H-bridge inverter. ") between the
# "...Fig. 1." sentence and "The below model is used ..." sentence.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("as represented in Fig. 1. The below model is used as a standard of comparison")
if ($found1) {
    $r1.Text = "as represented in Fig. 1. As can be observed the conventional cycloconverter using uses 2 separate converters called the P-converter and the N-converter; each performing like an H-bridge inverter. The below model is used as a standard of comparison"
}

# ---------------------------------------------------------------------------
# Change 2: "Fig, 2" / "Fig, 3" -> "Fig. 2" / "Fig. 3" (comma to period)
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Fig, 2", $true, $false, $false, $false, $false, $true, 1, $false, "Fig. 2", 2) | Out-Null

$r3 = $d.Content
$r3.Find.Execute("Fig, 3", $true, $false, $false, $false, $false, $true, 1, $false, "Fig. 3", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: remove "As is seen in the graphs..." sentence (keep the now-empty
# paragraph), remove the following blank paragraph entirely, and move the
# "_GoBack" bookmark from in front of "Before you begin..." to in front of
# "The below model is used..." (mirroring where the author's last edit
# landed).
# ---------------------------------------------------------------------------

# 3a. Locate the "As is seen..." paragraph and wipe its text only, leaving
#     the paragraph mark (and so the paragraph itself) intact.
$seenIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*As is seen in the graphs*") {
        $seenIndex = $i
        break
    }
}
if ($seenIndex -gt 0) {
    $p = $d.Paragraphs($seenIndex)
    $textOnly = $d.Range($p.Range.Start, $p.Range.End - 1)
    $textOnly.Delete()
}

# 3b. The paragraph right after it is an empty paragraph - remove it outright
#     by deleting its paragraph mark, which merges it (nothing) into the
#     following paragraph.
if ($seenIndex -gt 0) {
    $blank = $d.Paragraphs($seenIndex + 1)
    $blankMark = $d.Range($blank.Range.Start, $blank.Range.End)
    $blankMark.Delete()
}

# 3c. Move the "_GoBack" bookmark.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$r4 = $d.Content
$found4 = $r4.Find.Execute("The below model is used as a standard of comparison for the SiC model")
if ($found4) {
    $spot = $d.Range($r4.Start, $r4.Start)
    $d.Bookmarks.Add("_GoBack", $spot) | Out-Null
}
